# Turn the literal run "koji je na sjednici Skupštine ŠŠS-a," into the
# merge-field placeholder "{{ obrazl }}" followed by its own comma run,
# so that clicking primijeni_btn populates obrazl_textbox correctly.
$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("koji je na sjednici Skupštine ŠŠS-a,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found -and $rng.Find.Found) {
    $start = $rng.Start

    # Replace the whole run's text (keeps the existing "s2" character style).
    $rng.Text = "{{ obrazl }},"

    # Re-point to the just-replaced span and carve out everything except the
    # trailing comma; restyling that sub-range forces Word to split the run
    # into two runs, matching the target: "{{ obrazl }}" (apple-converted-space)
    # followed by "," (s2).
    $placeholderEnd = $rng.End - 1
    $sub = $d.Range($start, $placeholderEnd)
    $sub.Style = "apple-converted-space"
}
